$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play Artemis vs Medusa Free Slot
#    Game by Quickspin").
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs.Item(2)
$p2.Style = "Normal"

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our Artemis vs Medusa review to learn about Quickspin' + [char]39 + 's latest Greek mythology-inspired online slot game and play for free on desktop or mobile.</w:t></w:r></w:p>'
$p2.Range.InsertXML($metaXml) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Artemis vs Medusa Free Slot
#    Game by Quickspin" paragraph near the end of the document (the
#    Heading1 at the very top has the same text, so match on the
#    exact paragraph text *and* skip the Heading 1 styled one).
# ------------------------------------------------------------------
$dupTitle = "Play Artemis vs Medusa Free Slot Game by Quickspin"
$dupIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -ne "Heading 1") {
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $dupTitle) {
            $dupIndex = $i
            break
        }
    }
}
if ($dupIndex -ge 1) {
    $d.Paragraphs.Item($dupIndex).Range.Delete() | Out-Null
}

# ------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new
#    image-prompt copy, keeping its italic run formatting intact.
#    Scope the Find to that specific paragraph so the similarly
#    worded "Meta description" sentence above is left untouched.
# ------------------------------------------------------------------
$oldCopy = "Read our Artemis vs Medusa review to learn about Quickspin's latest Greek mythology-inspired online slot game and play for free on desktop or mobile."
$newCopy = "Create a cartoon-style feature image for the Artemis vs Medusa slot game that features a happy Maya warrior with glasses. The image should convey the epic battle between Artemis and Medusa with the Maya warrior standing confidently in the middle of the two opposing forces. The warrior should be wearing traditional Maya garb with a bow and arrow in hand, ready to join the battle. The background should depict a Greek temple in ruins with the sun setting behind it, giving the image an overall ancient and mythical feel. Make sure to add bold colors and details to capture attention and excitement."

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq $oldCopy) {
        $p.Range.Find.Execute($oldCopy, $true, $false, $false, $false, $false,
                               $true, 1, $false, $newCopy, 2) | Out-Null
        break
    }
}
